$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 632, shifting existing rows 632:663 down to 633:664
$ws.Rows.Item(632).Insert()

$newRow = 632

$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 45267
$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = 100112009
$ws.Cells.Item($newRow, 7).Value = "Acelga"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 230
$ws.Cells.Item($newRow, 11).Value = 3000
$ws.Cells.Item($newRow, 12).Value = 3500
$ws.Cells.Item($newRow, 13).Value = 3261
$ws.Cells.Item($newRow, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item($newRow, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($newRow, 16).Value = 544
$ws.Cells.Item($newRow, 17).Value = 6
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"

# Match the date cell formatting used by column D elsewhere (style index 2)
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
